# Edit: 
#  1) Swap the two theme parts' colour schemes (the "Integral"/Red-Violet
#     palette on the slide master's theme becomes the "Office"/blue palette
#     that used to live only on the notes master's theme, and vice versa).
#     In the live object model the only colour scheme PowerPoint lets an
#     automation client repaint is the one that is actually in effect for
#     the slides (Slide.ThemeColorScheme) -- there is no supported
#     Master/NotesMaster.ThemeColorScheme in the PowerPoint object model,
#     so we drive the swap from a normal slide, which edits the theme part
#     used by the slide master/slides.
#  2) Re-point every table on the deck from the old table style GUID to the
#     new one via Table.ApplyStyle.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Theme colours: replace the "Red Violet" accent palette that is
#    currently applied to the slides with the stock "Office" palette
#    (dk1/lt1 -- black/white -- are already identical in both palettes).
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72

# ---------------------------------------------------------------------
# 2) Tables: point every table at the new table-style GUID.
# ---------------------------------------------------------------------
$oldStyle = "{33102CE8-E1CD-4E54-88AA-DACC7953CE21}"
$newStyle = "{734FD098-3BE6-4426-A6BD-BEAB801FAD10}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyle)
        }
    }
}
